$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

# --- Row 3 ---
$ws.Range("A3").Value = 0
$ws.Range("A3").Font.Color = 255
$ws.Range("B3").Value = "Portcullis  - Not implemented"
$ws.Range("B3").Font.Color = 255

# --- Row 4 ---
$ws.Range("A4").Value = 1
$ws.Range("A4").Font.Color = 255
$ws.Range("B4").Value = "Cheval De Frise "

# --- Row 5 ---
$ws.Range("A5").Value = 2
$ws.Range("A5").Font.Color = 255
$ws.Range("B5").Value = "Moat "

# --- Row 6 ---
$ws.Range("A6").Value = 3
$ws.Range("A6").Font.Color = 255
$ws.Range("B6").Value = "Ramparts "

# --- Row 7 ---
$ws.Range("A7").Value = 4
$ws.Range("A7").Font.Color = 255
$ws.Range("B7").Value = "Drawbridge - Not implemented"
$ws.Range("B7").Font.Color = 255

# --- Row 8 ---
$ws.Range("A8").Value = 5
$ws.Range("A8").Font.Color = 255
$ws.Range("B8").Value = "Sally Port  - Not implemented"
$ws.Range("B8").Font.Color = 255

# --- Row 9 ---
$ws.Range("A9").Value = 6
$ws.Range("A9").Font.Color = 255
$ws.Range("B9").Value = "Rock Wall  - Not implemented"
$ws.Range("B9").Font.Color = 255

# --- New rows 10-14: copy formatting from row 3 (col A) / row 4 (col B) ---
$ws.Range("A3").Copy($ws.Range("A10"))
$ws.Range("B4").Copy($ws.Range("B10"))
$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Rough Terrain"

$ws.Range("A3").Copy($ws.Range("A11"))
$ws.Range("B4").Copy($ws.Range("B11"))
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "Corner Shot"

$ws.Range("A3").Copy($ws.Range("A12"))
$ws.Range("B4").Copy($ws.Range("B12"))
$ws.Range("A12").Value = 9
$ws.Range("B12").Value = "Low Bar (Close with /targetTrack)"

$ws.Range("A3").Copy($ws.Range("A13"))
$ws.Range("B4").Copy($ws.Range("B13"))
$ws.Range("A13").Value = 21
$ws.Range("B13").Value = "Low Bar (Far with /targetTrack)"

$ws.Range("A3").Copy($ws.Range("A14"))
$ws.Range("B4").Copy($ws.Range("B14"))
$ws.Range("A14").Value = "default"
$ws.Range("B14").Value = "Corner Shot"

# --- Restore default row heights on header rows (matches diff removing ht="14.45") ---
$ws.Range("A1:B2").EntireRow.AutoFit()

# --- Selection state ---
$ws.Range("A1:B14").Select()
